$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# E1/I1/J1/L1 text is unchanged; F1/G1/H1/K1 got new content (columns rotated)
$ws.Range("E1").Value = "Client"
$ws.Range("F1").Value = "Lob"
$ws.Range("G1").Value = "Process"
$ws.Range("H1").Value = "Product Name"
$ws.Range("K1").Value = "Municipality"
$ws.Range("M1").Value = "Tier"

# --- Row 2 ---
$ws.Range("A2").Value = 45436.041666666664
$ws.Range("B2").Value = "WFG18-001"
$ws.Range("C2").Value = "SIPL5316"
$ws.Range("D2").Value = "SIPL5688"
$ws.Range("E2").Value = "WFG Title"
$ws.Range("F2").Value = "Title"
$ws.Range("G2").Value = "Search"
$ws.Range("H2").Value = "Full Search"
$ws.Range("I2").Value = "AL"
$ws.Range("J2").Value = "Shelby"
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = "WIP"
$ws.Range("M2").Value = "Search(T1)"

# --- Row 3 ---
$ws.Range("A3").Value = 45439.083333333336
$ws.Range("B3").Value = "WFG18-002"
$ws.Range("C3").Value = "SIPL6118"
$ws.Range("D3").Value = "SIPL4167"
$ws.Range("E3").Value = "WFG Title"
$ws.Range("F3").Value = "Title"
$ws.Range("G3").Value = "Search"
$ws.Range("H3").Value = "Current Owner Search"
$ws.Range("I3").Value = "FL"
$ws.Range("J3").Value = "Clay"
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = "WIP"
$ws.Range("M3").Value = "Search(T2)"

# --- Cell formatting: black font on OrderID/Emp ID/Client columns for the data rows,
#     plus drop the left border on the Emp ID column (D) ---
$ws.Range("C2").Font.Color = 0
$ws.Range("E2").Font.Color = 0
$ws.Range("C3").Font.Color = 0
$ws.Range("E3").Font.Color = 0

$ws.Range("D2").Font.Color = 0
$ws.Range("D2").Borders.Item(7).LineStyle = -4142
$ws.Range("D3").Font.Color = 0
$ws.Range("D3").Borders.Item(7).LineStyle = -4142

# --- Column widths (content-driven resize of columns C and H) ---
$ws.Range("C1").EntireColumn.ColumnWidth = 19.6
$ws.Range("H1").EntireColumn.ColumnWidth = 18.1

# --- Selection state ---
$ws.Range("H5").Select()
